$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from column I (rows 4-14) to the new column J
$ws.Range("I4:I14").Copy()
$ws.Range("J4:J14").PasteSpecial(-4122)

# Fill in the values for the new column J (year 2021)
$ws.Range("J4").Value = 2021
$ws.Range("J5").Value = 1.5
$ws.Range("J6").Value = 0.3
$ws.Range("J7").Value = 0.8
$ws.Range("J8").Value = 0.6
$ws.Range("J9").Value = 1.8
$ws.Range("J10").Value = 0.5
$ws.Range("J11").Value = 0.8
$ws.Range("J12").Value = 1.9
$ws.Range("J13").Value = 4.4000000000000004
$ws.Range("J14").Value = 0.4

# Update the active cell selection to match the saved state
[void]$ws.Range("L10").Select()
